# SM-II (H) Attendance Sheet - add Day 5 (column K) attendance marks
# and hide the now-completed Section/Total/Day1-4 columns (D:J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 7-82 hold one participant each. Column K is "Day 5" (header row 6,
# K6 = 5). Every participant gets a new "P" (Present) or "A" (Absent) mark
# for that day. The rows below were marked Absent for day 5; everyone else
# was marked Present.
$absentRows = @(10, 13, 15, 16, 19, 24, 27, 34, 49, 65, 76)

for ($r = 7; $r -le 82; $r++) {
    if ($absentRows -contains $r) {
        $ws.Cells.Item($r, 11).Value = "A"
    } else {
        $ws.Cells.Item($r, 11).Value = "P"
    }
}

# Match the formatting already used for the other "day" columns (column I,
# which carries style index 29 for every row) so the new K column entries
# look identical to the rest of the attendance grid instead of the bare
# "empty cell" style that previously occupied K7:K82.
$ws.Range("I7:I82").Copy()
$ws.Range("K7:K82").PasteSpecial(-4122)

# The Section / Total Absence / Total Present / Day1-4 columns (D:J) are now
# "done" for this update, so hide them and keep the freshly-entered Day 5
# column (K) and everything after it visible.
$ws.Range("D1:J1").EntireColumn.Hidden = $true
